$wb = $excel.ActiveWorkbook

$amsin = $wb.Worksheets.Item("AMSIN")
$ams   = $wb.Worksheets.Item("AMS")

# ---------------------------------------------------------------------------
# Helper style template cells living on AMSIN (sheet "AMSIN"), already present
# in the workbook before any edits:
#   AMSIN!A2   -> plain/general data style (dates/text columns A,C,D,E,F,G)
#   AMSIN!B9   -> datetime style used for row-8-style timestamps
#   AMSIN!B10  -> datetime style used for the "normal" per-row timestamps
# We use Range.Copy(destination) to clone formatting onto new / changed
# cells, then overwrite Value2 with the real payload so the stored style
# index matches an existing one instead of minting a brand-new duplicate.
# ---------------------------------------------------------------------------

function Set-TextDateCell($range, [string]$text) {
    # Writing a literal "yyyy-mm-dd"-looking string straight into Value2
    # gets auto-parsed into a real date serial. Going through a text
    # formula and then collapsing it to a static value with PasteSpecial
    # (values only) keeps the cell as plain text instead.
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# ===========================================================================
# Sheet AMSIN ("AMSIN"): fix row 13, append rows 14 and 15
# ===========================================================================

# --- Row 13: correct the stored run-time value (tiny float fix) -----------
$amsin.Range("B13").Value2 = 44362.74926832176

# --- Row 14: new entry -----------------------------------------------------
Set-TextDateCell $amsin.Range("A14") "2021-06-17"

$amsin.Range("B10").Copy($amsin.Range("B14"))
$amsin.Range("B14").Value2 = 44364.73979796296

$amsin.Range("C14").Value2 = "whatsappbulk"
$amsin.Range("D14").Value2 = 77
$amsin.Range("E14").Value2 = 75
$amsin.Range("F14").Value2 = 2
$amsin.Range("G14").Value2 = 1.98

# --- Row 15: new entry -----------------------------------------------------
Set-TextDateCell $amsin.Range("A15") "2021-06-17"

$amsin.Range("B10").Copy($amsin.Range("B15"))
$amsin.Range("B15").Value2 = 44364.75321611918

$amsin.Range("C15").Value2 = "qwerty123rgrsn"
$amsin.Range("D15").Value2 = 77
$amsin.Range("E15").Value2 = 72
$amsin.Range("F15").Value2 = 5
$amsin.Range("G15").Value2 = 2.64

# ===========================================================================
# Sheet AMS: fill in row 8 (previously just an empty placeholder cell),
# append rows 9 and 10
# ===========================================================================

# --- Row 8 ------------------------------------------------------------------
Set-TextDateCell $ams.Range("A8") "2021-06-17"

$ams.Range("B8").Value2 = 44364.54754861111
$ams.Range("B8").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ams.Range("C8").Value2 = "145_data_hstry"
$ams.Range("D8").Value2 = 69
$ams.Range("E8").Value2 = 68
$ams.Range("F8").Value2 = 1
$ams.Range("G8").Value2 = 2.45

# --- Row 9 -------------------------------------------------------------------
Set-TextDateCell $ams.Range("A9") "2021-06-17"

$amsin.Range("B10").Copy($ams.Range("B9"))
$ams.Range("B9").Value2 = 44364.65389450231

$ams.Range("C9").Value2 = "145_hstry_data"
$ams.Range("D9").Value2 = 70
$ams.Range("E9").Value2 = 69
$ams.Range("F9").Value2 = 1
$ams.Range("G9").Value2 = 1.98

# --- Row 10 --------------------------------------------------------------
Set-TextDateCell $ams.Range("A10") "2021-06-17"

$amsin.Range("B10").Copy($ams.Range("B10"))
$ams.Range("B10").Value2 = 44364.72706061343

$ams.Range("C10").Value2 = "145_livetest"
$ams.Range("D10").Value2 = 70
$ams.Range("E10").Value2 = 63
$ams.Range("F10").Value2 = 7
$ams.Range("G10").Value2 = 2.62
